$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 70: Consecrating Congregation | Holy Water
$ws.Range("H70").Value = 3228.2856
$ws.Range("I70").Value = 1800
$ws.Range("J70").Value = 3466.3333
$ws.Range("K70").Value = 5400
$ws.Range("L70").Value = 10398.9999
$ws.Range("M70").Value = -5130
$ws.Range("N70").Value = -10938.9999

# ALC row 73: Curbing the Contagion (L) | Holy Water
$ws.Range("H73").Value = 3228.2856
$ws.Range("I73").Value = 1800
$ws.Range("J73").Value = 3466.3333
$ws.Range("K73").Value = 5400
$ws.Range("L73").Value = 10398.9999
$ws.Range("M73").Value = -4464
$ws.Range("N73").Value = -12270.9999

# ALC row 92: Whinier than the Sword | Enchanted Koppranickel Ink
$ws.Range("H92").Value = 1373.15
$ws.Range("J92").Value = 2014
$ws.Range("L92").Value = 2014
$ws.Range("N92").Value = -4510

# ALC row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 6041.1177
$ws.Range("I116").Value = 5464.923
$ws.Range("K116").Value = 5464.923
$ws.Range("M116").Value = -2022.923

# ALC row 135: For Tired Minds | Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 2144.4062
$ws.Range("I135").Value = 2137.8262
$ws.Range("K135").Value = 19240.4358
$ws.Range("M135").Value = -16705.4358

# ALC row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3458.532
$ws.Range("I138").Value = 2777.8333
$ws.Range("J138").Value = 3691.9143
$ws.Range("K138").Value = 8333.499899999999
$ws.Range("L138").Value = 11075.7429
$ws.Range("M138").Value = -3193.499899999999
$ws.Range("N138").Value = -21355.7429

$ws = $wb.Worksheets.Item("ARM")
# ARM row 110: Scheduled Maintenance | Deepgold Ingot
$ws.Range("H110").Value = 2280.2856
$ws.Range("I110").Value = 2160.3333
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 2160.3333
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = -115.3332999999998
$ws.Range("N110").Value = -7090

# ARM row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 7381.0557
$ws.Range("I132").Value = 3136.2273
$ws.Range("J132").Value = 14051.5
$ws.Range("K132").Value = 9408.6819
$ws.Range("L132").Value = 42154.5
$ws.Range("M132").Value = -6878.6819
$ws.Range("N132").Value = -47214.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 2416.8333
$ws.Range("I86").Value = 2521.4
$ws.Range("J86").Value = 1894
$ws.Range("K86").Value = 2521.4
$ws.Range("L86").Value = 1894
$ws.Range("M86").Value = -1398.4
$ws.Range("N86").Value = -4140

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 2416.8333
$ws.Range("I89").Value = 2521.4
$ws.Range("J89").Value = 1894
$ws.Range("K89").Value = 12607
$ws.Range("L89").Value = 9470
$ws.Range("M89").Value = -6991
$ws.Range("N89").Value = -20702

# BSM row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 46771.477
$ws.Range("I134").Value = 3415.7144
$ws.Range("J134").Value = 502007
$ws.Range("K134").Value = 10247.1432
$ws.Range("L134").Value = 1506021
$ws.Range("M134").Value = -7712.143199999999
$ws.Range("N134").Value = -1511091

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 841808.4
$ws.Range("I31").Value = 13855.5
$ws.Range("K31").Value = 13855.5
$ws.Range("M31").Value = -13560.5

# CRP row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 841808.4
$ws.Range("I34").Value = 13855.5
$ws.Range("K34").Value = 13855.5
$ws.Range("M34").Value = -13653.5

# CRP row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 4077.2307
$ws.Range("I132").Value = 4077.2307
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12231.6921
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9701.6921
$ws.Range("N132").ClearContents()

# CRP row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 386998.7
$ws.Range("I134").Value = 456094.88
$ws.Range("J134").Value = 6969.75
$ws.Range("K134").Value = 1368284.64
$ws.Range("L134").Value = 20909.25
$ws.Range("M134").Value = -1365749.64
$ws.Range("N134").Value = -25979.25

$ws = $wb.Worksheets.Item("CUL")
# CUL row 39: Bloody Good Tart, This | Blood Currant Tart
$ws.Range("H39").Value = 138243.62
$ws.Range("J39").Value = 355333
$ws.Range("L39").Value = 1065999
$ws.Range("N39").Value = -1066587

# CUL row 92: Oh No Udon | Gyr Abanian Flour
$ws.Range("H92").Value = 668020
$ws.Range("J92").Value = 792.3077
$ws.Range("L92").Value = 2376.9231
$ws.Range("N92").Value = -4872.9231

# CUL row 98: Sweet Kiss of Death | Rice Vinegar
$ws.Range("H98").Value = 699
$ws.Range("J98").Value = 598
$ws.Range("L98").Value = 1794
$ws.Range("N98").Value = -4790

# CUL row 107: Slippery Service | Frantoio Oil
$ws.Range("H107").Value = 664.1875
$ws.Range("I107").Value = 455.22223
$ws.Range("J107").Value = 932.8570999999999
$ws.Range("K107").Value = 1365.66669
$ws.Range("L107").Value = 2798.5713
$ws.Range("M107").Value = 554.33331
$ws.Range("N107").Value = -6638.5713

# CUL row 128: A Historical Flavor | Skyr
$ws.Range("H128").Value = 444662
$ws.Range("I128").Value = 444662
$ws.Range("K128").Value = 1333986
$ws.Range("M128").Value = -1329006

# CUL row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 2237.4614
$ws.Range("J132").Value = 2057.3845
$ws.Range("L132").Value = 18516.4605
$ws.Range("N132").Value = -23576.4605

$ws = $wb.Worksheets.Item("GSM")
# GSM row 43: Get the Green Stuff | Malachite Earrings
$ws.Range("H43").Value = 8635.625
$ws.Range("I43").Value = 4155
$ws.Range("K43").Value = 4155
$ws.Range("M43").Value = -4004

# GSM row 124: The Sage's Successor | Pewter Pendulums
$ws.Range("H124").Value = 75000
$ws.Range("J124").Value = 75000
$ws.Range("L124").Value = 75000
$ws.Range("N124").Value = -84820

# GSM row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 4762.636
$ws.Range("I126").Value = 4477.8
$ws.Range("K126").Value = 13433.4
$ws.Range("M126").Value = -10963.4

# GSM row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 32260444
$ws.Range("I132").Value = 33335738
$ws.Range("J132").Value = 1610
$ws.Range("K132").Value = 100007214
$ws.Range("L132").Value = 4830
$ws.Range("M132").Value = -100004684
$ws.Range("N132").Value = -9890

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16: Saddle Sore | Hard Leather
$ws.Range("H16").Value = 900.087
$ws.Range("J16").Value = 811.5
$ws.Range("L16").Value = 811.5
$ws.Range("N16").Value = -1151.5

# LTW row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 3142.7856
$ws.Range("I22").Value = 2444.3333
$ws.Range("J22").Value = 4400
$ws.Range("K22").Value = 2444.3333
$ws.Range("L22").Value = 4400
$ws.Range("M22").Value = -2149.3333
$ws.Range("N22").Value = -4990

# LTW row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 3142.7856
$ws.Range("I27").Value = 2444.3333
$ws.Range("J27").Value = 4400
$ws.Range("K27").Value = 2444.3333
$ws.Range("L27").Value = 4400
$ws.Range("M27").Value = -2337.3333
$ws.Range("N27").Value = -4614

# LTW row 68: You Could Say It's a Moving Target | Wyvern Leather
$ws.Range("H68").Value = 3110
$ws.Range("I68").Value = 2954.5454
$ws.Range("J68").Value = 3537.5
$ws.Range("K68").Value = 2954.5454
$ws.Range("L68").Value = 3537.5
$ws.Range("M68").Value = -2205.5454
$ws.Range("N68").Value = -5035.5

# LTW row 71: They Call It Bloody Mary (L) | Wyvern Leather
$ws.Range("H71").Value = 3110
$ws.Range("I71").Value = 2954.5454
$ws.Range("J71").Value = 3537.5
$ws.Range("K71").Value = 14772.727
$ws.Range("L71").Value = 17687.5
$ws.Range("M71").Value = -11028.727
$ws.Range("N71").Value = -25175.5

# LTW row 82: Trainin' the Neck | Dragon Leather
$ws.Range("H82").Value = 1831.75
$ws.Range("I82").Value = 1774.6666
$ws.Range("J82").Value = 2003
$ws.Range("K82").Value = 1774.6666
$ws.Range("L82").Value = 2003
$ws.Range("M82").Value = -1413.6666
$ws.Range("N82").Value = -2725

# LTW row 85: Training Is Only Skintight (L) | Dragon Leather
$ws.Range("H85").Value = 1831.75
$ws.Range("I85").Value = 1774.6666
$ws.Range("J85").Value = 2003
$ws.Range("K85").Value = 1774.6666
$ws.Range("L85").Value = 2003
$ws.Range("M85").Value = -526.6666
$ws.Range("N85").Value = -4499

# LTW row 101: A Stitch in Time | Marid Leather Gloves of Healing
$ws.Range("H101").Value = 49997.2
$ws.Range("J101").Value = 49997.2
$ws.Range("L101").Value = 49997.2
$ws.Range("N101").Value = -56487.2

# LTW row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 3934477.8
$ws.Range("I132").Value = 460253.5
$ws.Range("K132").Value = 1380760.5
$ws.Range("M132").Value = -1378230.5

$ws = $wb.Worksheets.Item("WVR")
# WVR row 103: To the Tops | Serge Gambison of Healing
$ws.Range("H103").Value = 101854.664
$ws.Range("J103").Value = 101854.664
$ws.Range("L103").Value = 101854.664
$ws.Range("N103").Value = -104198.664

# WVR row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 297735.12
$ws.Range("I132").Value = 2423.48
$ws.Range("J132").Value = 1118045.2
$ws.Range("K132").Value = 7270.440000000001
$ws.Range("L132").Value = 3354135.6
$ws.Range("M132").Value = -4740.440000000001
$ws.Range("N132").Value = -3359195.6

# WVR row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 3825.7144
$ws.Range("J136").Value = 6364.143
$ws.Range("L136").Value = 19092.429
$ws.Range("N136").Value = -24192.429
